# Updated TestData for Portugal Market
#
# 1) Duplicate the "Czech" sheet (same A1:D21 template as every other market
#    sheet) and drop the copy right after "Swiss" so it lands at the end of
#    the tab strip.
# 2) Rename the copy to "Portugal", fill in its market name / NGC ticket
#    cells, and select B2 (so the new sheet matches the freshly-edited look
#    the author left it in).
# 3) Germany's selection is reset to the full used range (no more "last
#    active cell" pointer).
# 4) The active tab moves to the new last sheet ("Swiss" loses
#    tabSelected, "Portugal" gains it) - this happens automatically because
#    Copy() activates the newly created sheet.

$wb = $excel.ActiveWorkbook

# Germany's selection is updated first - Select() activates the sheet it's
# called on, and the final active sheet needs to be the new "Portugal" tab,
# so this has to happen before the new sheet is created/activated below.
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D21").Select() | Out-Null

$czech = $wb.Worksheets.Item("Czech")
$swiss = $wb.Worksheets.Item("Swiss")

$czech.Copy($null, $swiss)

$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Set B4 before B2 so new shared strings land in the same order as the
# authored workbook (NGC ticket first, market name second).
$portugal.Range("B4").Value = "NGC-3479/T2434/T2437"
$portugal.Range("B2").Value = "Portugal Market"

$portugal.Range("B2").Select() | Out-Null
